# Updates cryptos list values (price / volume change %) to reflect the
# latest scrape, per commit "Updated cryptos list ... with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing it to be stored as text, so that
# numeric-looking strings (e.g. "159.35", "1.00", "0.0905") are not
# reinterpreted by Excel as numbers (which would introduce floating point
# artifacts and change the cell type). Resetting the style back to
# "Normal" afterwards avoids leaving a stray quote-prefix style applied.
function Set-TextValue($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Simple price (column D) / change% (column E) updates, keyed by row number.
# A value of $null means that column is left unchanged for that row.
$updates = @(
    @{Row=2;  D="69.348.56";  E="  -0.04%  "},
    @{Row=3;  D="3.690.66";   E="  -0.02%  "},
    @{Row=4;  D=$null;        E="  +0.03%  "},
    @{Row=5;  D="681.13";     E="  -1.32%  "},
    @{Row=6;  D="159.35";     E="  -1.73%  "},
    @{Row=7;  D=$null;        E="  -0.08%  "},
    @{Row=8;  D=$null;        E="  -0.57%  "},
    @{Row=10; D="7.13";       E="  -3.22%  "},
    @{Row=11; D=$null;        E="  -0.21%  "},
    @{Row=12; D=$null;        E="  -2.42%  "},
    @{Row=13; D="4.312.66";   E=$null},
    @{Row=14; D="32.48";      E="  -2.17%  "},
    @{Row=15; D="3.682.08";   E="  -0.25%  "},
    @{Row=16; D="69.343.33";  E="  -0.10%  "},
    @{Row=17; D=$null;        E="  +2.05%  "},
    @{Row=18; D="16.02";      E="  -0.60%  "},
    @{Row=19; D=$null;        E="  -0.80%  "},
    @{Row=20; D="468.85";     E="  -1.87%  "},
    @{Row=21; D="9.92";       E="  -0.46%  "},
    @{Row=22; D="0.655";      E="  -0.83%  "},
    @{Row=23; D="79.92";      E="  -0.01%  "},
    @{Row=24; D="3.836.83";   E="  +0.00%  "},
    @{Row=25; D=$null;        E="  -0.04%  "},
    @{Row=26; D=$null;        E="  -4.39%  "},
    @{Row=27; D="10.91";      E="  -3.35%  "},
    @{Row=28; D="9.13";       E="  -3.04%  "},
    @{Row=29; D=$null;        E="  -0.93%  "},
    @{Row=30; D=$null;        E="  -3.93%  "},
    @{Row=31; D="6.60";       E="  -3.19%  "},
    @{Row=32; D=$null;        E="  -3.06%  "},
    @{Row=33; D=$null;        E="  +0.14%  "},
    @{Row=34; D="26.94";      E="  +0.09%  "},
    @{Row=35; D="3.678.35";   E="  +0.48%  "},
    @{Row=36; D=$null;        E="  -6.86%  "},
    @{Row=37; D="8.32";       E="  -1.26%  "},
    @{Row=38; D=$null;        E="  -0.41%  "},
    @{Row=42; D="0.0905";     E="  -1.74%  "},
    @{Row=43; D="171.19";     E="  +4.37%  "},
    @{Row=44; D="0.943";      E=$null},
    @{Row=45; D="47.52";      E="  -1.15%  "},
    @{Row=46; D=$null;        E="  -4.70%  "},
    @{Row=47; D=$null;        E="  -1.49%  "},
    @{Row=48; D=$null;        E="  -2.10%  "},
    @{Row=49; D=$null;        E="  -1.95%  "},
    @{Row=50; D=$null;        E="  -1.85%  "},
    @{Row=51; D="7.82";       E="  -2.64%  "}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        Set-TextValue $ws.Cells.Item($r, 4) $u.D
    }
    if ($null -ne $u.E) {
        Set-TextValue $ws.Cells.Item($r, 5) $u.E
    }
}

# Rows 40 and 41 swapped coins (Stacks <-> FirstDigitalUSD) along with
# their refreshed price / change values.
Set-TextValue $ws.Cells.Item(40, 2) "FirstDigitalUSD"
Set-TextValue $ws.Cells.Item(40, 3) "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Cells.Item(40, 4) "1.00"
Set-TextValue $ws.Cells.Item(40, 5) "  -0.13%  "

Set-TextValue $ws.Cells.Item(41, 2) "Stacks"
Set-TextValue $ws.Cells.Item(41, 3) "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Cells.Item(41, 4) "2.23"
Set-TextValue $ws.Cells.Item(41, 5) "  -4.45%  "
